$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header / group info ---
$ws.Range("D2").Value = 10

$ws.Range("D4").Value = 2201793
$ws.Range("G4").Value = "Tiago José Figueira Pires Rodrigues dos Reis"

$ws.Range("D6").Value = 2201790
$ws.Range("G6").Value = "Daniel Marques Gonçalves"

$ws.Range("D7").Value = 2182185
$ws.Range("G7").Value = "João Pedro Da Rocha Valverde Martins"

# --- Funcionamento do Grupo scores (rows 26-28) ---
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 4
$ws.Range("L26").Value = 3

$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 4
$ws.Range("L27").Value = 4

$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 2
$ws.Range("L28").Value = 1

# --- Grupo de Funcionalidades do Projeto (rows 31-38) ---
$ws.Range("I31").Value = "Completo"
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = 4
$ws.Range("L31").Value = 4

$ws.Range("I32").Value = "Completo"
$ws.Range("J32").Value = 4
$ws.Range("K32").Value = 3
$ws.Range("L32").Value = 2

$ws.Range("I33").Value = "Parcial"

$ws.Range("I35").Value = "Parcial"

$ws.Range("I37").Value = "Completo"
$ws.Range("J37").Value = 1
$ws.Range("K37").Value = 4
$ws.Range("L37").Value = 4
